$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06108550355870301
$ws.Range("H2").Value = -5.002467573744219
$ws.Range("I2").Value = -4.389989509649465
$ws.Range("G3").Value = 0.06188942395223233
$ws.Range("H3").Value = 10.10973737703844
$ws.Range("G4").Value = -0.01896204520231863
$ws.Range("H4").Value = 6.644293452391715
$ws.Range("G5").Value = -0.008884569031835427
$ws.Range("H5").Value = 22.2787874868036
$ws.Range("G6").Value = -0.01795737899452797
$ws.Range("H6").Value = -60.42601937104326
$ws.Range("G7").Value = -0.009858991866272688
$ws.Range("H7").Value = -74.07223512492286
$ws.Range("G8").Value = -0.01065108269748892
$ws.Range("H8").Value = -87.21704234440222
$ws.Range("G9").Value = -0.00550713131900317
$ws.Range("H9").Value = -0.2545914657072099
$ws.Range("G10").Value = -0.06116487087508717
$ws.Range("H10").Value = 2.946417134835258
$ws.Range("G11").Value = -0.06699325044482141
$ws.Range("H11").Value = -4.549455778698159
$ws.Range("G12").Value = -0.3925890600288665
$ws.Range("H12").Value = 0.5104722400808851
$ws.Range("G13").Value = -0.4049963313122755
$ws.Range("H13").Value = -3.283983841473785
$ws.Range("G14").Value = -0.0398388850309683
$ws.Range("H14").Value = -386.1801790361301
$ws.Range("G15").Value = -0.01205533760053082
$ws.Range("H15").Value = 73.39874239792266
$ws.Range("G16").Value = 0.1355250622491215
$ws.Range("H16").Value = -0.8954955960497976
$ws.Range("G17").Value = 0.1505074642154252
$ws.Range("H17").Value = 7.913368589149361
$ws.Range("G18").Value = 0.1263188586979595
$ws.Range("H18").Value = 7.351731648834058
$ws.Range("G19").Value = 0.1279803735820699
$ws.Range("H19").Value = -0.5350771746361078
$ws.Range("G20").Value = 0.09025591333353654
$ws.Range("H20").Value = 1.7144158851114
$ws.Range("G21").Value = 0.08767414448794696
$ws.Range("H21").Value = 0.6695514477513889
$ws.Range("G22").Value = -0.09406813257922748
$ws.Range("H22").Value = -0.6257416989547001
$ws.Range("G23").Value = -0.1046261089369834
$ws.Range("H23").Value = -3.135778248637576
$ws.Range("G24").Value = 0.1589331470242742
$ws.Range("H24").Value = -1.339561886837346
$ws.Range("G25").Value = 0.1705825295573944
$ws.Range("H25").Value = -0.007739001108871615
$ws.Range("G26").Value = 0.085840162414438
$ws.Range("H26").Value = -5.309908137725635
$ws.Range("G27").Value = 0.08451360963456629
$ws.Range("H27").Value = -1.684629160876963
$ws.Range("G28").Value = -0.1397098143927103
$ws.Range("H28").Value = -1.503466631101927
$ws.Range("G29").Value = -0.1374400556852353
$ws.Range("H29").Value = 1.687921131763558
$ws.Range("G30").Value = 0.05156709513859661
$ws.Range("H30").Value = -0.8616304624179104
$ws.Range("G31").Value = 0.05124122925274138
$ws.Range("H31").Value = 16.93730568991531
$ws.Range("G32").Value = 0.117289067406854
$ws.Range("H32").Value = 7.897723302226564
$ws.Range("G33").Value = 0.1172584743721433
$ws.Range("H33").Value = -5.509252962520649
$ws.Range("G34").Value = -0.01219393778703231
$ws.Range("H34").Value = 21.91579036470122
$ws.Range("G35").Value = -0.01427233760097591
$ws.Range("H35").Value = 14.72730078019862
$ws.Range("G36").Value = 0.03691548951543461
$ws.Range("H36").Value = 0.4041133603752248
$ws.Range("G37").Value = 0.04122427056832208
$ws.Range("H37").Value = 15.52058696876169
$ws.Range("G38").Value = 0.09861961925524829
$ws.Range("H38").Value = -1.677540100709837
$ws.Range("G39").Value = 0.1095194060495277
$ws.Range("H39").Value = 12.44290474156998
$ws.Range("G40").Value = 0.02675170681104557
$ws.Range("H40").Value = -20.58890769993582
$ws.Range("G41").Value = 0.03100219187421486
$ws.Range("H41").Value = -3.775857545081227
$ws.Range("G42").Value = 0.1186482841517387
$ws.Range("H42").Value = -1.866869275115442
$ws.Range("G43").Value = 0.1293947311683842
$ws.Range("H43").Value = 1.259364030344056
$ws.Range("G44").Value = 0.03518764850840481
$ws.Range("H44").Value = -11.29158102710069
$ws.Range("G45").Value = 0.03135694861267752
$ws.Range("H45").Value = 0.624457600803317
$ws.Range("G46").Value = 0.05699436226728728
$ws.Range("H46").Value = 0.6653352031202864
$ws.Range("G47").Value = 0.06277282266496996
$ws.Range("H47").Value = 6.99270998730823
$ws.Range("G48").Value = 0.04280702334554923
$ws.Range("H48").Value = -13.09510752991362
$ws.Range("G49").Value = 0.04932204671811011
$ws.Range("H49").Value = 8.238940174502416
$ws.Range("G50").Value = 0.02790690510487296
$ws.Range("H50").Value = 5.362815859381839
$ws.Range("G51").Value = 0.02352371427149231
$ws.Range("H51").Value = -16.03432067408623
$ws.Range("G52").Value = -0.08599072996590115
$ws.Range("H52").Value = 1.073253643432666
$ws.Range("G53").Value = -0.08084153278083568
$ws.Range("H53").Value = -0.7795567936724608
$ws.Range("G54").Value = 0.04767648865149864
$ws.Range("H54").Value = -4.697917682514405
$ws.Range("G55").Value = 0.0494222770045336
$ws.Range("H55").Value = -12.1928983746002
$ws.Range("G56").Value = 0.0492426560675418
$ws.Range("H56").Value = -0.3828275648560823
$ws.Range("G57").Value = 0.03915236443594951
$ws.Range("H57").Value = 3.075998470091529
$ws.Range("G58").Value = 0.0487183509940922
$ws.Range("H58").Value = -15.440159497035
$ws.Range("G59").Value = 0.06288526004518866
$ws.Range("H59").Value = 10.29926316692926
$ws.Range("G60").Value = 0.02169850612478036
$ws.Range("H60").Value = -20.98140384087836
$ws.Range("G61").Value = 0.02728983509047989
$ws.Range("H61").Value = 2.2174456781833
$ws.Range("G62").Value = 0.0665901919504634
$ws.Range("H62").Value = 6.634338455353355
$ws.Range("G63").Value = 0.06630204041330832
$ws.Range("H63").Value = 3.787210339795307
$ws.Range("G64").Value = 0.02374466316035677
$ws.Range("H64").Value = -14.4017433319012
$ws.Range("G65").Value = 0.03431252204347111
$ws.Range("H65").Value = -3.146532143849884
$ws.Range("G66").Value = 0.07370284307011246
$ws.Range("H66").Value = -5.123809344581868
$ws.Range("G67").Value = 0.08439601008840614
$ws.Range("H67").Value = 7.005348766465308
$ws.Range("G68").Value = -0.01784624477603438
$ws.Range("H68").Value = 17.93204753329445
$ws.Range("G69").Value = -0.01665821328705916
$ws.Range("H69").Value = 12.97245917025791
$ws.Range("G70").Value = 0.07058759711031934
$ws.Range("H70").Value = -1.942417113893392
$ws.Range("G71").Value = 0.07241727961440149
$ws.Range("H71").Value = -8.81851763967504
$ws.Range("G72").Value = -0.1499125893100363
$ws.Range("H72").Value = 2.427009285453315
$ws.Range("G73").Value = -0.1420974321188938
$ws.Range("H73").Value = 7.171702058114819
$ws.Range("G74").Value = 0.154069244349664
$ws.Range("H74").Value = 2.426585587712979
$ws.Range("G75").Value = 0.1512483879362196
$ws.Range("H75").Value = 0.5324543987060277
$ws.Range("G76").Value = -0.001767336318942978
$ws.Range("H76").Value = -70.51632373416177
$ws.Range("G77").Value = 0.001564786624721301
$ws.Range("H77").Value = 170.8725613386662
$ws.Range("G78").Value = 0.09266024926201422
$ws.Range("H78").Value = 3.00361426877494
$ws.Range("G79").Value = 0.09491513019667207
$ws.Range("H79").Value = -2.049748766856716
$ws.Range("G80").Value = -0.2217780217350827
$ws.Range("H80").Value = -2.469422308735426
$ws.Range("G81").Value = -0.2159300991202397
$ws.Range("H81").Value = -1.321459543000836
$ws.Range("G82").Value = 0.1714835460651409
$ws.Range("H82").Value = 2.306366567380052
$ws.Range("G83").Value = 0.1724087587243887
$ws.Range("H83").Value = -2.056578663164524
$ws.Range("G84").Value = 0.1066920210508802
$ws.Range("H84").Value = 0.5298460765620355
$ws.Range("G85").Value = 0.1106277610172464
$ws.Range("H85").Value = 5.804345221689731
